$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header rows: shift the "12-month ending" year labels forward by one year
# (drop 1396/12, columns now show 1397/12 .. 1401/12)
$ws.Cells.Item(8, 5).Value = "دوازده ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 6).Value = "دوازده ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 7).Value = "دوازده ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 8).Value = "دوازده ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 9).Value = "دوازده ماهه منتهی به 1401/12"

$ws.Cells.Item(28, 5).Value = "دوازده ماهه منتهی به 1397/12"
$ws.Cells.Item(28, 6).Value = "دوازده ماهه منتهی به 1398/12"
$ws.Cells.Item(28, 7).Value = "دوازده ماهه منتهی به 1399/12"
$ws.Cells.Item(28, 8).Value = "دوازده ماهه منتهی به 1400/12"
$ws.Cells.Item(28, 9).Value = "دوازده ماهه منتهی به 1401/12"

$ws.Cells.Item(48, 5).Value = "دوازده ماهه منتهی به 1397/12"
$ws.Cells.Item(48, 6).Value = "دوازده ماهه منتهی به 1398/12"
$ws.Cells.Item(48, 7).Value = "دوازده ماهه منتهی به 1399/12"
$ws.Cells.Item(48, 8).Value = "دوازده ماهه منتهی به 1400/12"
$ws.Cells.Item(48, 9).Value = "دوازده ماهه منتهی به 1401/12"

$ws.Cells.Item(61, 5).Value = "دوازده ماهه منتهی به 1397/12"
$ws.Cells.Item(61, 6).Value = "دوازده ماهه منتهی به 1398/12"
$ws.Cells.Item(61, 7).Value = "دوازده ماهه منتهی به 1399/12"
$ws.Cells.Item(61, 8).Value = "دوازده ماهه منتهی به 1400/12"
$ws.Cells.Item(61, 9).Value = "دوازده ماهه منتهی به 1401/12"

$ws.Cells.Item(81, 5).Value = "دوازده ماهه منتهی به 1397/12"
$ws.Cells.Item(81, 6).Value = "دوازده ماهه منتهی به 1398/12"
$ws.Cells.Item(81, 7).Value = "دوازده ماهه منتهی به 1399/12"
$ws.Cells.Item(81, 8).Value = "دوازده ماهه منتهی به 1400/12"
$ws.Cells.Item(81, 9).Value = "دوازده ماهه منتهی به 1401/12"

# Unit label change: "لیتر" (liter) -> "بطری/عدد" (bottle/piece)
$ws.Cells.Item(13, 3).Value = "بطری/عدد"
$ws.Cells.Item(18, 3).Value = "بطری/عدد"
$ws.Cells.Item(53, 3).Value = "بطری/عدد / ریال"
$ws.Cells.Item(57, 3).Value = "بطری/عدد / ریال"

# Data rows: shift values one column to the left (E<-F, F<-G, G<-H, H<-I)
# and populate the new column I with the new 1401/12 figures
$rowsData = @{
    11 = @(1118769, 1422659, 1890157, 1542266, 1764571)
    12 = @(201328, 166058, 350730, 452707, 216990)
    13 = @(774491, 939156, 954832, 79206290, 66006118)
    14 = @(2094588, 2527873, 3195719, 81201263, 67987679)
    16 = @(23080, 2990, 19980, 30130, 0)
    17 = @(24702, 30522, 14387, 19296, 23522)
    18 = @("-", 4332, 9240, 19465900, 0)
    19 = @(47782, 37844, 43607, 19515326, 23522)
    22 = @(0, 0, 0, "-", "-")
    23 = @(0, 0, 0, 0, 0)
    24 = @(2142370, 2565717, 3239326, 100716589, 68011201)
    31 = @(875296, 1469393, 4024340, 5850601, 9273090)
    32 = @(415998, 793653, 849944, 1553945, 682776)
    33 = @(209254, 375619, 672588, 811343, 1514676)
    34 = @(1500548, 2638665, 5546872, 8215889, 11470542)
    36 = @(6237, 4117, 25818, 69300, 0)
    37 = @(583517, 890538, 1554659, 1180976, 1613998)
    38 = @("-", 21519, 163112, 100761, 0)
    39 = @(589754, 916174, 1743589, 1351037, 1613998)
    42 = @(0, 0, 0, "-", "-")
    43 = @(0, 0, 0, 0, 0)
    44 = @(2090302, 3554839, 7290461, 9566926, 13084540)
    51 = @(782374, 1032850, 2129104, 3793510, 5255153)
    52 = @(2066270, 4779372, 2423357, 3432562, 3146578)
    53 = @(270183, 399954, 704405, 10243, 22948)
    55 = @(270234, 1376923, 1292192, 2300033, 0)
    56 = @(23622257, 29176922, 108059985, 61203151, 68616529)
    57 = @("-", 4967452, 17652814, 5176, 0)
    64 = @(-730796, -1006415, -2194922, -2942913, -5076138)
    65 = @(-189877, -334767, -383125, -622191, -542893)
    66 = @(-84319, -215366, -322517, -478052, -882503)
    67 = @(-1004992, -1556548, -2900564, -4043156, -6501534)
    69 = @(-13934, -2427, -15821, -44291, 0)
    70 = @(-186408, -404390, -270700, -268718, -467304)
    71 = @("-", -11381, -50528, -53203, 0)
    72 = @(-200342, -418198, -337049, -366212, -467304)
    75 = @(0, 0, 0, "-", "-")
    76 = @(0, 0, 0, 0, 0)
    77 = @(-1205334, -1974746, -3237613, -4409368, -6968838)
    84 = @(144500, 462978, 1829418, 2907688, 4196952)
    85 = @(226121, 458886, 466819, 931754, 139883)
    86 = @(124935, 160253, 350071, 333291, 632173)
    87 = @(495556, 1082117, 2646308, 4172733, 4969008)
    89 = @(-7697, 1690, 9997, 25009, 0)
    90 = @(397109, 486148, 1283959, 912258, 1146694)
    91 = @("-", 10138, 112584, 47558, 0)
    92 = @(389412, 497976, 1406540, 984825, 1146694)
    95 = @(884968, 1580093, 4052848, 5157558, 6115702)
}

foreach ($r in $rowsData.Keys) {
    $vals = $rowsData[$r]
    $col = 5
    foreach ($v in $vals) {
        $ws.Cells.Item($r, $col).Value = $v
        $col = $col + 1
    }
}